{"js": "// Update the \"payment hours\" paragraph:\n//  - Mon-Sat 8:30am-5:30pm / 3 days notice for Sunday appointments\n//    becomes\n//  - Mon-Fri 9am-5pm, Sat 9am-2pm / at least 3 days notice for Sunday appointments\nconst oldText =\n  \"LOS PAGOS DEBER\u00c1N REALIZARSE DE LUNES A S\u00c1BADO, ENTRE LAS 8:30 A.M. Y LAS 5:30 P.M. PARA EFECTUAR UN PAGO EN DOMINGO, SER\u00c1 INDISPENSABLE PROGRAMAR UNA CITA CON TRES D\u00cdAS DE ANTICIPACI\u00d3N. CADA PAGO DEBER\u00c1 SER NOTIFICADO Y CONFIRMADO AL N\u00daMERO TELEF\u00d3NICO 951 189 9298.\";\nconst newText =\n  \"LOS PAGOS DEBER\u00c1N REALIZARSE DE LUNES A VIERNES, EN UN HORARIO DE 9:00 A.M. A 5:00 P.M., Y EN S\u00c1BADO DE 9:00 A. M. A 2:00 P. M. PARA EFECTUAR UN PAGO EN DOMINGO, SER\u00c1 INDISPENSABLE PROGRAMAR UNA CITA CON AL MENOS TRES D\u00cdAS DE ANTICIPACI\u00d3N. CADA PAGO DEBER\u00c1 SER NOTIFICADO Y CONFIRMADO AL N\u00daMERO TELEF\u00d3NICO 951 189 9298.\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Update the \"payment hours\" paragraph:\n#  - Mon-Sat 8:30am-5:30pm / 3 days notice for Sunday appointments\n#    becomes\n#  - Mon-Fri 9am-5pm, Sat 9am-2pm / at least 3 days notice for Sunday appointments\n$d = $word.ActiveDocument\n\n$oldText = \"LOS PAGOS DEBER\u00c1N REALIZARSE DE LUNES A S\u00c1BADO, ENTRE LAS 8:30 A.M. Y LAS 5:30 P.M. PARA EFECTUAR UN PAGO EN DOMINGO, SER\u00c1 INDISPENSABLE PROGRAMAR UNA CITA CON TRES D\u00cdAS DE ANTICIPACI\u00d3N. CADA PAGO DEBER\u00c1 SER NOTIFICADO Y CONFIRMADO AL N\u00daMERO TELEF\u00d3NICO 951 189 9298.\"\n$newText = \"LOS PAGOS DEBER\u00c1N REALIZARSE DE LUNES A VIERNES, EN UN HORARIO DE 9:00 A.M. A 5:00 P.M., Y EN S\u00c1BADO DE 9:00 A. M. A 2:00 P. M. PARA EFECTUAR UN PAGO EN DOMINGO, SER\u00c1 INDISPENSABLE PROGRAMAR UNA CITA CON AL MENOS TRES D\u00cdAS DE ANTICIPACI\u00d3N. CADA PAGO DEBER\u00c1 SER NOTIFICADO Y CONFIRMADO AL N\u00daMERO TELEF\u00d3NICO 951 189 9298.\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 2) | Out-Null\n"}
